$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated bus voltage magnitude (vm_pu) results for rows 2-25 (case with 380 kV done).
# Columns B:F and I:N are updated per row; column G (slack bus, value 1) and the
# blank column H are left unchanged.

$bfValues = @(
    @(1.02,1.046622029432734,1.066276933474887,1.056001738619828,1.071115961139446),
    @(1.02,1.047523065515943,1.066977362202182,1.056783674533134,1.071934318636481),
    @(1.02,1.048106725103795,1.067430900236477,1.057290507400888,1.072464622561623),
    @(1.02,1.048352244794476,1.067621640984599,1.057503786066545,1.072687745314112),
    @(1.02,1.048393477339952,1.067653671445833,1.057539608531465,1.072725219283503),
    @(1.02,1.048110005165387,1.067433448637664,1.057293356433029,1.072467603222044),
    @(1.02,1.046926407726861,1.06651358060314,1.056265816776087,1.071392367969957),
    @(1.02,1.044845636297485,1.06489513713743,1.054461883215228,1.069503658031989),
    @(1.02,1.043461817197436,1.063817948478855,1.053263880901409,1.068248654556102),
    @(1.02,1.042863422690441,1.063351956739746,1.05274624957266,1.067706228009754),
    @(1.02,1.042641275112487,1.063178934255801,1.052554146951727,1.067504898441528),
    @(1.02,1.042688920944061,1.063216045078468,1.052595345954614,1.067548077418609),
    @(1.02,1.042845057389748,1.063337653258003,1.052730366858789,1.067689582937542),
    @(1.02,1.042941274517764,1.063412589081445,1.052813580099601,1.067776789328022),
    @(1.02,1.04350154781493,1.063848884197491,1.053298257952732,1.068284674807404),
    @(1.02,1.043853209805945,1.064122679136935,1.053602582366502,1.068603526498722),
    @(1.02,1.044058406193755,1.064282421239692,1.053780196785566,1.068789603488896),
    @(1.02,1.044128386031166,1.064336896263189,1.053840776833109,1.06885306719428),
    @(1.02,1.043815471735293,1.064093299158784,1.053569920143015,1.068569306766912),
    @(1.02,1.042799075697761,1.063301840799838,1.052690601897165,1.067647908883975),
    @(1.02,1.042160738237942,1.062804611488867,1.052138716515184,1.067069468652383),
    @(1.02,1.042499065150529,1.063068164405227,1.052431188232001,1.067376026816516),
    @(1.02,1.04383252370445,1.064106574576256,1.053584678465051,1.068584768893748),
    @(1.02,1.045382979045832,1.065313238644726,1.054927436294413,1.069991213631526)
)

$inValues = @(
    @(1.053250765733669,1.051675444115624,1.068987915750333,1.058740647635617,1.073813992218905,1.053168943875127),
    @(1.053514238606107,1.052224861482665,1.069504193011164,1.059336216902908,1.07444883673121,1.053719141477849),
    @(1.053683688648075,1.052580342785495,1.069837894455299,1.059721796679053,1.074859739762261,1.054075127604903),
    @(1.053754676968539,1.052729779137434,1.06997809392408,1.059883942564952,1.075032509721689,1.054224776173603),
    @(1.053766581625572,1.052754869656033,1.070001628795502,1.059911170357693,1.075061520058538,1.054249902323616),
    @(1.053684638174844,1.052582339594197,1.069839768157384,1.059723963092785,1.074862048221343,1.054077127249303),
    @(1.053340021342234,1.051861127477939,1.069162468712715,1.058941879478416,1.074028515943168,1.053354890929115),
    @(1.052724875061296,1.050590093306738,1.067966249466966,1.057565389191405,1.072560681166137,1.052082051743602),
    @(1.052309522633538,1.049742694939214,1.067167014842608,1.05664891079679,1.071582855635197,1.051233449973197),
    @(1.052128434129268,1.049375766640174,1.066820538175535,1.056252362328961,1.071159638365319,1.050866000593885),
    @(1.052060984564512,1.049239474409284,1.066691782157368,1.056105111869819,1.071002466307113,1.050729514812395),
    @(1.052075461094542,1.049268709479068,1.066719403415504,1.056136695522135,1.071036178880161,1.050758791399332),
    @(1.0521228625013,1.049364500649099,1.066809896361837,1.056240189628375,1.071146645859893,1.050854718603812),
    @(1.052152043551123,1.049423520954339,1.066865644254472,1.056303961804881,1.071214712189348,1.050913822724655),
    @(1.052321514861779,1.049767046877417,1.067190000987376,1.056675234692602,1.071610947251989,1.051257836493946),
    @(1.052427488802511,1.049982532641254,1.067393354514866,1.056908203668375,1.071859546325614,1.051473628272284),
    @(1.052489182091525,1.050108221839645,1.067511928231431,1.057044118685623,1.072004567963753,1.051599495963754),
    @(1.052510197612547,1.050151078572633,1.067552352155673,1.057090466929616,1.072054019603226,1.051642413558219),
    @(1.052416131158177,1.049959413049958,1.067371540625425,1.056883205380314,1.071832872124691,1.051450475848517),
    @(1.052108909075407,1.049336292492498,1.066783250066732,1.056209711932562,1.071114115256418,1.050826470388393),
    @(1.05191467505888,1.04894451982838,1.066413026850991,1.055786522605444,1.070662376739892,1.050434141362155),
    @(1.052017743469076,1.049152204788287,1.066609321107727,1.056010837929904,1.070901835030824,1.050642121258526),
    @(1.052421263558143,1.049969859798229,1.067381397503622,1.056894500943809,1.071844924992456,1.051460937432368),
    @(1.052884834863391,1.050918698649173,1.068275816105589,1.057921042373558,1.072940029277278,1.052411123743314)
)

$startRow = 2
$bfCols = @(2,3,4,5,6)     # columns B..F
$inCols = @(9,10,11,12,13,14)  # columns I..N

for ($i = 0; $i -lt $bfValues.Count; $i++) {
    $r = $startRow + $i
    $bfRow = $bfValues[$i]
    for ($j = 0; $j -lt $bfCols.Count; $j++) {
        $ws.Cells.Item($r, $bfCols[$j]).Value = $bfRow[$j]
    }
    $inRow = $inValues[$i]
    for ($j = 0; $j -lt $inCols.Count; $j++) {
        $ws.Cells.Item($r, $inCols[$j]).Value = $inRow[$j]
    }
}

Write-Output "Updated B2:F25 and I2:N25 with new vm_pu results"
